# dictionary.xlsx: merge dateStart/dateEnd/timeStart/timeEnd terms into
# simple start/end terms (process overview without date/time split, in
# line with "Process overview and drivers without summation").
#
# Old rows (A = calculation label, B = bill label):
#   1 label used for calculation | label on bill   (header)
#   2 driver                     | Fahrer
#   3 car                        | Fahrzeug
#   4 dateStart                  | Datum Beginn
#   5 dateEnd                    | Datum Ende
#   6 timeStart                  | Zeit Beginn
#   7 timeEnd                    | Zeit Ende
#   8 duration                   | Dauer
#   9 distance                   | Distanz
#  10 cost                       | Kosten
#
# New rows:
#   1 label used for calculation | label on bill   (header)
#   2 driver                     | Fahrer
#   3 car                        | Fahrzeug
#   4 start                      | Beginn
#   5 end                        | Ende
#   6 duration                   | Dauer
#   7 distance                   | Distanz
#   8 cost                       | Kosten

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the dateStart/dateEnd/timeStart/timeEnd dictionary rows as the
# simplified start/end rows, shifting duration/distance/cost up two rows.
$ws.Cells.Item(4, 1).Value2 = "start"
$ws.Cells.Item(4, 2).Value2 = "Beginn"
$ws.Cells.Item(5, 1).Value2 = "end"
$ws.Cells.Item(5, 2).Value2 = "Ende"
$ws.Cells.Item(6, 1).Value2 = "duration"
$ws.Cells.Item(6, 2).Value2 = "Dauer"
$ws.Cells.Item(7, 1).Value2 = "distance"
$ws.Cells.Item(7, 2).Value2 = "Distanz"
$ws.Cells.Item(8, 1).Value2 = "cost"
$ws.Cells.Item(8, 2).Value2 = "Kosten"

# The old distance/cost rows (9 and 10) are now redundant since their
# content moved up to rows 7 and 8 above - remove the two trailing rows
# so the table ends at row 8 (dimension A1:B10 -> A1:B8).
$ws.Rows.Item(9).Delete() | Out-Null
$ws.Rows.Item(9).Delete() | Out-Null

# Match the new active cell/selection recorded in the sheet view.
$ws.Range("A6").Select() | Out-Null

# Slightly narrower columns (minor relayout that accompanied the edit).
$ws.Columns.Item(1).ColumnWidth = 22.1166666666667
$ws.Columns.Item(2).ColumnWidth = 10.5085034013605

# Reflect the updated sheet-tab area ratio from the source workbook view.
$win = $wb.Windows.Item(1)
$win.TabRatio = 0.985 | Out-Null
